$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-13 (row 4 unchanged) — columns D (Fecha), J (Volumen),
# K (Precio mínimo), L (Precio máximo), M (Precio promedio ponderado), P (Precio $/Kg)
$rows = @(
    @{ Row=2;  D=44229; J=120; K=44000; L=45000; M=44500; P=3423 },
    @{ Row=3;  D=44397; J=140; K=12500; L=13000; M=12750; P=981  },
    @{ Row=5;  D=44580; J=160; K=11000; L=12000; M=11500; P=885  },
    @{ Row=6;  D=44389; J=120; K=12000; L=13000; M=12500; P=962  },
    @{ Row=7;  D=44592; J=120; K=12000; L=13000; M=12500; P=962  },
    @{ Row=8;  D=44469; J=140; K=13000; L=14000; M=13500; P=1038 },
    @{ Row=9;  D=44616; J=120; K=19000; L=20000; M=19500; P=1500 },
    @{ Row=10; D=44320; J=160; K=19000; L=20000; M=19500; P=1500 },
    @{ Row=11; D=44764; J=200; K=12000; L=13000; M=12500; P=962  },
    @{ Row=12; D=44159; J=100; K=23000; L=24000; M=23500; P=1808 },
    @{ Row=13; D=44379; J=120; K=12000; L=13000; M=12667; P=974  }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("D$n").Value = $r.D
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("P$n").Value = $r.P
}
